$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2022" data column (S), copying formatting from the adjacent
# 2021 column (R) so the new cells match the existing header/data styling.
$ws.Range("R3").Copy($ws.Range("S3"))
$ws.Range("S3").Value = 2022

$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 0.071025550219041236

# Columns A:C get a uniform width.
$ws.Range("A1:C1").ColumnWidth = 32.71

# Move the active selection.
$ws.Range("F14").Select() | Out-Null
